# Swap the order of "System" and the email address in the "Recorded By"
# column (G) wherever both appear together as "System, dnasr281@gmail.com".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$oldText = "System, dnasr281@gmail.com"
$newText = "dnasr281@gmail.com, System"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G ("Recorded By")
    if ($cell.Value2 -eq $oldText) {
        $cell.Value2 = $newText
    }
}
